$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused rows (35-39) from the bottom of the table
$ws.Rows("35:39").Delete()

# Rewrite cluster names and active-case counts for rows 2-34 with the refreshed data
$clusters = @(
    @('3323 Villa Maria Catholic Homes St Bernadette''s Aged Care Sunshine North', 13),
    @('3376 Royal Freemasons Coppin Centre Melbourne', 11),
    @('3601 Baptcare Westhaven community', 29),
    @('3653 Fronditha Thalpori St Albans Aged Care', 40),
    @('44121 Wallan Primary School Wallan', 15),
    @('44165 Greenvale Primary School', 21),
    @('44234 Lucknow Primary School Bairnsdale', 18),
    @('44321 Maiden Gully Primary School Maiden Gully', 11),
    @('44395 Buln Buln Primary School', 11),
    @('44507 Inverloch Primary School', 11),
    @('44811 Dandenong North Primary School Dandenong', 20),
    @('44853 St Albans North Primary School', 12),
    @('45158 Rowellyn Park Primary School Carrum Downs', 17),
    @('45161 Woodlands Primary School Langwarrin', 10),
    @('45249 Creekside K-9 College Caroline Springs', 15),
    @('45585 Mount Ridley College Craigieburn', 11),
    @('45695 Sacred Heart Primary School Yarrawonga', 30),
    @('4574 Village Glen Aged Care Residences Mornington', 21),
    @('45809 St Finbar''s Primary School Brighton East', 12),
    @('45967 St Clement of Rome School Bulleen', 10),
    @('46037 Nazareth Catholic Primary School Grovedale', 16),
    @('46050 Our Lady''s Catholic Primary School Craigieburn', 34),
    @('46125 Our Lady of the Southern Cross Primary School Manor Lakes', 26),
    @('46322 Minaret College Officer Campus Officer', 13),
    @('46328 Ilim College Dallas Primary Campus Inverloch', 12),
    @('46390 Al Siraat College Epping', 47),
    @('52380 Al Iman College Melton South', 27),
    @('52786 Hume Anglican Grammar Donnybrook Campus', 16),
    @('Hampton Park Primary School Hampton Park', 10),
    @('Ilim College Dallas Main Campus Dallas Oct', 27),
    @('Islamic College of Melbourne Tarneit Oct Nov', 15),
    @('Melton Willows Melton', 10),
    @('Morwell Park Primary School Morwell', 10)
)

for ($i = 0; $i -lt $clusters.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $clusters[$i][0]
    $ws.Cells.Item($row, 2).Value = $clusters[$i][1]
}

Write-Output "Updated $($clusters.Length) cluster rows"
